$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 17: "parrallel to serial" level-shifter part (adcin2 config support)
$ws.Range("A17").Value = "parrallel to serial"
$ws.Range("B17").Value = "https://www.mouser.ee/ProductDetail/Texas-Instruments/SN74HCS16507PWR?qs=DPoM0jnrROUIT0uZVZ3kvw%3D%3D"
$ws.Range("C17").Value = 0.43
$ws.Range("D17").Value = 1

# Row 19: swap the USB-PD negotiation chip link to the in-stock DK variant
$ws.Range("B19").Value = "https://www.mouser.ee/ProductDetail/Texas-Instruments/TPS65988DKRSHR?qs=DPoM0jnrROUevQj%2FLwa4Vw%3D%3D"

# Move the active selection to the newly edited row
$ws.Range("B19").Select()
